# Regenerate the "K" (strike count) column (G) of the save_data sheet.
#
# The source pipeline now derives the K column from the option chain's
# actual strike count instead of the legacy "Strike#" field, and the
# std/mean + s_vals used to build it were recalculated. The per-row K
# values below are the freshly computed results that replace the old
# column G values (header "K" in G1, row 2.. data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row (1-based worksheet row) -> new K value
$kValues = @{
    2 = 1;
    3 = 1;
    4 = 3;
    5 = 1;
    6 = 0;
    7 = 1;
    8 = 2;
    9 = 2;
    10 = 1;
    11 = 2;
    12 = 2;
    13 = 3;
    14 = 1;
    15 = 1;
    16 = 1;
    17 = 0;
    18 = 2;
    19 = 2;
    20 = 0;
    21 = 2;
    22 = 0;
    23 = 3;
    24 = 1;
    25 = 2;
    26 = 0;
    27 = 0;
    28 = 0;
    30 = 0;
    31 = 1;
    32 = 0;
    33 = 1;
    34 = 3;
    35 = 1;
    36 = 1;
    37 = 1;
    38 = 2;
    39 = 1;
    40 = 2;
    41 = 1;
    42 = 1;
    43 = 2;
    44 = 2;
    45 = 3;
    46 = 3;
    47 = 0;
    48 = 1;
    49 = 1;
    50 = 3;
    51 = 1;
    52 = 0;
    53 = 1;
    54 = 0;
    55 = 3;
    56 = 2;
    57 = 0;
    58 = 1;
    59 = 1;
    60 = 1;
    61 = 3;
    62 = 0;
    63 = 1;
    65 = 2;
    66 = 1;
    68 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value2 = $kValues[$row]
}

Write-Host ("Updated " + $kValues.Count + " K (column G) values")
